$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (G2=5489)
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

# Row 34 (G34=2160)
$ws.Range("H34").Value = 14148.571
$ws.Range("I34").Value = 13173.333
$ws.Range("K34").Value = 13173.333
$ws.Range("M34").Value = -12970.333

# Row 36 (G36=2160)
$ws.Range("H36").Value = 14148.571
$ws.Range("I36").Value = 13173.333
$ws.Range("K36").Value = 13173.333
$ws.Range("M36").Value = -12458.333

# Row 53 (G53=5479)
$ws.Range("H53").Value = 456.4762
$ws.Range("I53").Value = 460
$ws.Range("K53").Value = 460
$ws.Range("M53").Value = 177

# Row 86 (G86=12603)
$ws.Range("H86").Value = 2514.4443
$ws.Range("I86").Value = 1522.1818
$ws.Range("J86").Value = 4073.7144
$ws.Range("K86").Value = 1522.1818
$ws.Range("L86").Value = 4073.7144
$ws.Range("M86").Value = -399.1818000000001
$ws.Range("N86").Value = -6319.7144

# Row 89 (G89=12603)
$ws.Range("H89").Value = 2514.4443
$ws.Range("I89").Value = 1522.1818
$ws.Range("J89").Value = 4073.7144
$ws.Range("K89").Value = 7610.909000000001
$ws.Range("L89").Value = 20368.572
$ws.Range("M89").Value = -1994.909000000001
$ws.Range("N89").Value = -31600.572

# Row 112 (G112=27960)
$ws.Range("H112").Value = 862.75
$ws.Range("J112").Value = 832.5161000000001
$ws.Range("L112").Value = 2497.5483
$ws.Range("N112").Value = -4713.5483

# Row 131 (G131=36108)
$ws.Range("H131").Value = 4252.2
$ws.Range("I131").Value = 3940.25
$ws.Range("K131").Value = 11820.75
$ws.Range("M131").Value = -6780.75

# Row 137 (G137=44013)
$ws.Range("H137").Value = 4176.8066
$ws.Range("I137").Value = 3539.1
$ws.Range("J137").Value = 5336.273
$ws.Range("K137").Value = 10617.3
$ws.Range("L137").Value = 16008.819
$ws.Range("M137").Value = -8067.299999999999
$ws.Range("N137").Value = -21108.819

# Row 138 (G138=44169)
$ws.Range("H138").Value = 2749.9
$ws.Range("I138").Value = 1092.0667
$ws.Range("J138").Value = 4407.7334
$ws.Range("K138").Value = 3276.2001
$ws.Range("L138").Value = 13223.2002
$ws.Range("M138").Value = 1863.7999
$ws.Range("N138").Value = -23503.2002

# Row 141 (G141=44161)
$ws.Range("H141").Value = 1030.5555
$ws.Range("I141").Value = 1030.5555
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 3091.6665
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2088.3335
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (G32=44147)
$ws.Range("H32").Value = 6412054.5
$ws.Range("I32").Value = 7813288.5
$ws.Range("J32").Value = 6412.143
$ws.Range("K32").Value = 7813288.5
$ws.Range("L32").Value = 6412.143
$ws.Range("M32").Value = -7813001.5
$ws.Range("N32").Value = -6986.143

# Row 61 (G61=43999)
$ws.Range("H61").Value = 3000.6428
$ws.Range("I61").Value = 2819.8667
$ws.Range("K61").Value = 2819.8667
$ws.Range("M61").Value = -2607.8667

# Row 132 (G132=43997)
$ws.Range("H132").Value = 2580.7551
$ws.Range("I132").Value = 2312.8838
$ws.Range("K132").Value = 6938.651400000001
$ws.Range("M132").Value = -4408.651400000001

# Row 136 (G136=43999)
$ws.Range("H136").Value = 3000.6428
$ws.Range("I136").Value = 2819.8667
$ws.Range("K136").Value = 8459.6001
$ws.Range("M136").Value = -5909.6001

$ws = $wb.Worksheets.Item("BSM")
# Row 31 (G31=1694)
$ws.Range("H31").Value = 15249.75
$ws.Range("I31").Value = 15499.5
$ws.Range("J31").Value = 15000
$ws.Range("K31").Value = 15499.5
$ws.Range("L31").Value = 15000
$ws.Range("M31").Value = -15247.5
$ws.Range("N31").Value = -15504

# Row 134 (G134=43998)
$ws.Range("H134").Value = 2189.4092
$ws.Range("I134").Value = 1157.8379
$ws.Range("K134").Value = 3473.5137
$ws.Range("M134").Value = -938.5137

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (G16=27691)
$ws.Range("H16").Value = 1626.5714
$ws.Range("I16").Value = 1402.625
$ws.Range("K16").Value = 1402.625
$ws.Range("M16").Value = -1115.625

# Row 58 (G58=44021)
$ws.Range("H58").Value = 1897.0465
$ws.Range("I58").Value = 1268.7587
$ws.Range("J58").Value = 3198.5
$ws.Range("K58").Value = 1268.7587
$ws.Range("L58").Value = 3198.5
$ws.Range("M58").Value = -1065.7587
$ws.Range("N58").Value = -3604.5

# Row 113 (G113=27691)
$ws.Range("H113").Value = 1626.5714
$ws.Range("I113").Value = 1402.625
$ws.Range("K113").Value = 1402.625
$ws.Range("M113").Value = 767.375

# Row 136 (G136=44021)
$ws.Range("H136").Value = 1897.0465
$ws.Range("I136").Value = 1268.7587
$ws.Range("J136").Value = 3198.5
$ws.Range("K136").Value = 3806.2761
$ws.Range("L136").Value = 9595.5
$ws.Range("M136").Value = -1256.2761
$ws.Range("N136").Value = -14695.5

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (G4=4650)
$ws.Range("H4").Value = 73263896
$ws.Range("I4").Value = 154405020
$ws.Range("J4").Value = 3714355.2
$ws.Range("K4").Value = 463215060
$ws.Range("L4").Value = 11143065.6
$ws.Range("M4").Value = -463214948
$ws.Range("N4").Value = -11143289.6

$ws = $wb.Worksheets.Item("GSM")
# Row 113 (G113=27710)
$ws.Range("H113").Value = 1799
$ws.Range("I113").Value = 1799
$ws.Range("K113").Value = 1799
$ws.Range("M113").Value = 371

# Row 132 (G132=44008)
$ws.Range("H132").Value = 24403942
$ws.Range("I132").Value = 37049790
$ws.Range("K132").Value = 111149370
$ws.Range("M132").Value = -111146840

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (G7=36249)
$ws.Range("H7").Value = 2807.35
$ws.Range("I7").Value = 1902.5264
$ws.Range("J7").Value = 19999
$ws.Range("K7").Value = 1902.5264
$ws.Range("L7").Value = 19999
$ws.Range("M7").Value = -1790.5264
$ws.Range("N7").Value = -20223

# Row 126 (G126=36249)
$ws.Range("H126").Value = 2807.35
$ws.Range("I126").Value = 1902.5264
$ws.Range("J126").Value = 19999
$ws.Range("K126").Value = 5707.5792
$ws.Range("L126").Value = 59997
$ws.Range("M126").Value = -3237.5792
$ws.Range("N126").Value = -64937

$ws = $wb.Worksheets.Item("WVR")
# Row 2 (G2=3307)
$ws.Range("H2").Value = 130953780
$ws.Range("I2").Value = 9522140
$ws.Range("J2").Value = 333339840
$ws.Range("K2").Value = 9522140
$ws.Range("L2").Value = 333339840
$ws.Range("M2").Value = -9522028
$ws.Range("N2").Value = -333340064

# Row 41 (G41=21725)
$ws.Range("H41").Value = 8077.6
$ws.Range("J41").Value = 9597
$ws.Range("L41").Value = 9597
$ws.Range("N41").Value = -10377

# Row 108 (G108=25661)
$ws.Range("H108").Value = 100593.164
$ws.Range("J108").Value = 100593.164
$ws.Range("L108").Value = 100593.164
$ws.Range("N108").Value = -108273.164

# Row 122 (G122=36208)
$ws.Range("H122").Value = 2401.4
$ws.Range("J122").Value = 4000
$ws.Range("L122").Value = 12000
$ws.Range("N122").Value = -16900

# Row 126 (G126=36210)
$ws.Range("H126").Value = 2307.4348
$ws.Range("I126").Value = 2366.6365
$ws.Range("J126").Value = 1005
$ws.Range("K126").Value = 7099.9095
$ws.Range("L126").Value = 3015
$ws.Range("M126").Value = -4629.9095
$ws.Range("N126").Value = -7955

# Row 132 (G132=44029)
$ws.Range("H132").Value = 1882.4
$ws.Range("I132").Value = 1789.0938
$ws.Range("J132").Value = 2255.625
$ws.Range("K132").Value = 5367.2814
$ws.Range("L132").Value = 6766.875
$ws.Range("M132").Value = -2837.2814
$ws.Range("N132").Value = -11826.875

# Row 133 (G133=41869)
$ws.Range("H133").Value = 66249.875
$ws.Range("J133").Value = 66249.875
$ws.Range("L133").Value = 66249.875
$ws.Range("N133").Value = -76369.875
